# Logged Week 16 and performed season sim from Week 17
# Updates cumulative season totals on the "Rushing" and "Receiving" sheets.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet -----------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# T.Lance (row 6): RZATT 4 -> 5
$rushing.Range("E2").Value = 13
$rushing.Range("E6").Value = 5

# J.Wilson (row 8): 1DATT 39->48, 2DATT 21->26, RZATT 6->9
$rushing.Range("C8").Value = 48
$rushing.Range("D8").Value = 26
$rushing.Range("F8").Value = 9

# D.Samuel (row 9): 1DATT 18->20, 2DATT 17->19, 3DATT 5->6, RZATT 10->11
$rushing.Range("C9").Value = 20
$rushing.Range("D9").Value = 19
$rushing.Range("E9").Value = 6
$rushing.Range("F9").Value = 11

# --- Receiving sheet -----------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# J.Hasty (row 4): Short Target 13->16, Short Comp 11->14
$receiving.Range("C4").Value = 16
$receiving.Range("D4").Value = 14

# K.Juszczyk (row 5): Short Target 23->26, Short Comp 21->23, Deep Target 3->5, Deep Comp 1->2
$receiving.Range("C5").Value = 26
$receiving.Range("D5").Value = 23
$receiving.Range("E5").Value = 5
$receiving.Range("F5").Value = 2

# J.Wilson (row 6): Short Target 5->8, Short Comp 4->7, RZ Target 2->3, RZ Comp 1->2
$receiving.Range("C6").Value = 8
$receiving.Range("D6").Value = 7
$receiving.Range("G6").Value = 3
$receiving.Range("H6").Value = 2

# D.Samuel (row 7): Short Target 73->81, Short Comp 42->48, Deep Target 25->27, Deep Comp 18->20
$receiving.Range("C7").Value = 81
$receiving.Range("D7").Value = 48
$receiving.Range("E7").Value = 27
$receiving.Range("F7").Value = 20

# B.Aiyuk (row 8): Short Target 47->52, Short Comp 31->35, RZ Target 8->9, RZ Comp 4->5
$receiving.Range("C8").Value = 52
$receiving.Range("D8").Value = 35
$receiving.Range("G8").Value = 9
$receiving.Range("H8").Value = 5

# J.Jennings (row 11): Short Target 21->24, Short Comp 12->14, Deep Target 4->6
$receiving.Range("C11").Value = 24
$receiving.Range("D11").Value = 14
$receiving.Range("E11").Value = 6

# G.Kittle (row 12): Short Target 80->83, Short Comp 66->68
$receiving.Range("C12").Value = 83
$receiving.Range("D12").Value = 68
